# Actualiza base de datos EC: el periodo de mora reportado para los
# trabajadores pasa de 2507 a 2508 (se agrega el nuevo periodo / "parte 1
# de nuevos estado de cuenta").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Columna E ("Periodo Mora") se almacena como texto, no como numero,
# para las tres filas de detalle de la tabla.
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"
